# PawnShop project update:
#  - Rename the original board sheet to InitBoard_Traditional
#  - Add a new InitBoard_PawnShop sheet holding just the two Kings'
#    starting squares (the minimal "pawn shop" variant setup)
#  - Bold the header row on both sheets
#  - Restore per-sheet selection / active-tab state

$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename the traditional board sheet -------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "InitBoard_Traditional"

# --- Sheet 2: new PawnShop init board, inserted right after sheet 1 ------
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "InitBoard_PawnShop"

$ws2.Range("A1").Value = "Side"
$ws2.Range("B1").Value = "Role"
$ws2.Range("C1").Value = "StartFile"
$ws2.Range("D1").Value = "StartRank"
$ws2.Range("A1:D1").Font.Bold = $true

# Only the kings start on the board in the PawnShop variant
$ws2.Range("A2").Value = "White"
$ws2.Range("B2").Value = "King"
$ws2.Range("C2").Value = 1
$ws2.Range("D2").Value = "e"

$ws2.Range("A3").Value = "Black"
$ws2.Range("B3").Value = "King"
$ws2.Range("C3").Value = 8
$ws2.Range("D3").Value = "e"

$ws2.Range("E13").Select() | Out-Null

# --- Back on sheet 1: bold header, restore selection/active tab ----------
$ws1.Activate() | Out-Null
$ws1.Range("A1:D1").Font.Bold = $true
$ws1.Range("I27").Select() | Out-Null
